$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 717.7778
$ws.Range("I80").Value = 337.63635
$ws.Range("J80").Value = 1315.1428
$ws.Range("K80").Value = 1012.90905
$ws.Range("L80").Value = 3945.4284
$ws.Range("M80").Value = -14.90904999999998
$ws.Range("N80").Value = -5941.428400000001

$ws.Range("H83").Value = 717.7778
$ws.Range("I83").Value = 337.63635
$ws.Range("J83").Value = 1315.1428
$ws.Range("K83").Value = 3038.72715
$ws.Range("L83").Value = 11836.2852
$ws.Range("M83").Value = 1953.27285
$ws.Range("N83").Value = -21820.2852

$ws.Range("H93").Value = 41831.58
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 41831.58
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 41831.58
$ws.Range("N93").Value = -46823.58

$ws.Range("H129").Value = 1221.6146
$ws.Range("I129").Value = 382.85715
$ws.Range("J129").Value = 1287.5842
$ws.Range("K129").Value = 1148.57145
$ws.Range("L129").Value = 3862.7526
$ws.Range("M129").Value = 3851.42855
$ws.Range("N129").Value = -13862.7526

$ws.Range("H132").Value = 24882474
$ws.Range("I132").Value = 25901704
$ws.Range("J132").Value = 5007503
$ws.Range("K132").Value = 77705112
$ws.Range("L132").Value = 15022509
$ws.Range("M132").Value = -77702582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4627.1763
$ws.Range("I32").Value = 4679.1816
$ws.Range("J32").Value = 4407.154
$ws.Range("K32").Value = 4679.1816
$ws.Range("L32").Value = 4407.154
$ws.Range("M32").Value = -4392.1816

$ws.Range("H61").Value = 2649.5454
$ws.Range("I61").Value = 2649.5454
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2649.5454
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2437.5454

$ws.Range("H119").Value = 32749.285
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 32749.285
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 32749.285
$ws.Range("N119").Value = -42425.285

$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800

$ws.Range("H132").Value = 2725.8333
$ws.Range("I132").Value = 2274.9092
$ws.Range("J132").Value = 4379.222
$ws.Range("K132").Value = 6824.7276
$ws.Range("L132").Value = 13137.666
$ws.Range("M132").Value = -4294.7276

$ws.Range("H136").Value = 2649.5454
$ws.Range("I136").Value = 2649.5454
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7948.6362
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5398.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2059.75
$ws.Range("I86").Value = 2040.875
$ws.Range("J86").Value = 2097.5
$ws.Range("K86").Value = 2040.875
$ws.Range("L86").Value = 2097.5
$ws.Range("M86").Value = -917.875
$ws.Range("N86").Value = -4343.5

$ws.Range("H89").Value = 2059.75
$ws.Range("I89").Value = 2040.875
$ws.Range("J89").Value = 2097.5
$ws.Range("K89").Value = 10204.375
$ws.Range("L89").Value = 10487.5
$ws.Range("M89").Value = -4588.375
$ws.Range("N89").Value = -21719.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 243691.86
$ws.Range("I31").Value = 1127283.8
$ws.Range("J31").Value = 2712.25
$ws.Range("K31").Value = 1127283.8
$ws.Range("L31").Value = 2712.25
$ws.Range("M31").Value = -1126988.8
$ws.Range("N31").Value = -3302.25

$ws.Range("H34").Value = 243691.86
$ws.Range("I34").Value = 1127283.8
$ws.Range("J34").Value = 2712.25
$ws.Range("K34").Value = 1127283.8
$ws.Range("L34").Value = 2712.25
$ws.Range("M34").Value = -1127081.8
$ws.Range("N34").Value = -3116.25

$ws.Range("H51").Value = 21271
$ws.Range("I51").Value = 7750
$ws.Range("J51").Value = 24651.25
$ws.Range("K51").Value = 7750
$ws.Range("L51").Value = 24651.25
$ws.Range("M51").Value = -7014
$ws.Range("N51").Value = -26123.25

$ws.Range("H61").Value = 21271
$ws.Range("I61").Value = 7750
$ws.Range("J61").Value = 24651.25
$ws.Range("K61").Value = 7750
$ws.Range("L61").Value = 24651.25
$ws.Range("M61").Value = -7402
$ws.Range("N61").Value = -25347.25

$ws.Range("H62").Value = 3666.3333
$ws.Range("I62").Value = 3499.625
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3499.625
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2875.625
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3666.3333
$ws.Range("I65").Value = 3499.625
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 17498.125
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -14378.125
$ws.Range("N65").Value = -31240

$ws.Range("H134").Value = 7201.778
$ws.Range("I134").Value = 7503.2666
$ws.Range("J134").Value = 5694.3335
$ws.Range("K134").Value = 22509.7998
$ws.Range("L134").Value = 17083.0005
$ws.Range("M134").Value = -19974.7998
$ws.Range("N134").Value = -22153.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2222.842
$ws.Range("I122").Value = 743.6667
$ws.Range("J122").Value = 2905.5386
$ws.Range("K122").Value = 6693.0003
$ws.Range("L122").Value = 26149.8474
$ws.Range("M122").Value = -4243.0003
$ws.Range("N122").Value = -31049.8474

$ws.Range("H131").Value = 822.8586
$ws.Range("I131").Value = 529
$ws.Range("J131").Value = 841.8172
$ws.Range("K131").Value = 1587
$ws.Range("L131").Value = 2525.4516
$ws.Range("M131").Value = 3453
$ws.Range("N131").Value = -12605.4516

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2990.5
$ws.Range("I80").Value = 2929.2856
$ws.Range("J80").Value = 3133.3333
$ws.Range("K80").Value = 2929.2856
$ws.Range("L80").Value = 3133.3333
$ws.Range("M80").Value = -1931.2856
$ws.Range("N80").Value = -5129.3333

$ws.Range("H83").Value = 2990.5
$ws.Range("I83").Value = 2929.2856
$ws.Range("J83").Value = 3133.3333
$ws.Range("K83").Value = 14646.428
$ws.Range("L83").Value = 15666.6665
$ws.Range("M83").Value = -9654.428
$ws.Range("N83").Value = -25650.6665

$ws.Range("H102").Value = 4667.6665
$ws.Range("I102").Value = 3501.5
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 3501.5
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = -1879.5
$ws.Range("N102").Value = -10244

$ws.Range("H132").Value = 2707.8
$ws.Range("I132").Value = 2221.2424
$ws.Range("J132").Value = 4045.8333
$ws.Range("K132").Value = 6663.7272
$ws.Range("L132").Value = 12137.4999
$ws.Range("M132").Value = -4133.7272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 962.2771
$ws.Range("I68").Value = 899.8642
$ws.Range("J68").Value = 3490
$ws.Range("K68").Value = 899.8642
$ws.Range("L68").Value = 3490
$ws.Range("M68").Value = -150.8642
$ws.Range("N68").Value = -4988

$ws.Range("H71").Value = 962.2771
$ws.Range("I71").Value = 899.8642
$ws.Range("J71").Value = 3490
$ws.Range("K71").Value = 4499.321
$ws.Range("L71").Value = 17450
$ws.Range("M71").Value = -755.3209999999999
$ws.Range("N71").Value = -24938

$ws.Range("H122").Value = 8601.799999999999
$ws.Range("I122").Value = 4004
$ws.Range("J122").Value = 9751.25
$ws.Range("K122").Value = 12012
$ws.Range("L122").Value = 29253.75
$ws.Range("M122").Value = -9562
$ws.Range("N122").Value = -34153.75

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 41980
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 41980
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 41980
$ws.Range("N125").Value = -51820

$ws.Range("H127").Value = 29294.285
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 29294.285
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 29294.285
$ws.Range("N127").Value = -39214.285

$ws.Range("H128").Value = 41997.145
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 41997.145
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 41997.145
$ws.Range("N128").Value = -51957.145

$ws.Range("H130").Value = 39957.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 39957.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 39957.5
$ws.Range("N130").Value = -49997.5

$ws.Range("H132").Value = 3712.2144
$ws.Range("I132").Value = 2666.9565
$ws.Range("J132").Value = 8520.4
$ws.Range("K132").Value = 8000.869499999999
$ws.Range("L132").Value = 25561.2
$ws.Range("M132").Value = -5470.869499999999

$ws.Range("H136").Value = 4668.269
$ws.Range("I136").Value = 2419.2856
$ws.Range("J136").Value = 7292.0835
$ws.Range("K136").Value = 7257.8568
$ws.Range("L136").Value = 21876.2505
$ws.Range("M136").Value = -4707.8568
$ws.Range("N136").Value = -26976.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3346.4583
$ws.Range("I132").Value = 2021.0667
$ws.Range("J132").Value = 5555.4443
$ws.Range("K132").Value = 6063.2001
$ws.Range("L132").Value = 16666.3329
$ws.Range("M132").Value = -3533.2001
$ws.Range("N132").Value = -21726.3329
